# Update handback/handoff timestamps to reflect regenerated report values.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2)
$wsOverview.Range("G2").Value = "2016-08-31 08:13:34"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 2
$wsZhCn.Range("H2").Value = "2016-08-31 08:13:23"
$wsZhCn.Range("K2").Value = "2016-08-31 08:14:49"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 2
$wsDeDe.Range("H2").Value = "2016-08-31 08:13:34"
$wsDeDe.Range("K2").Value = "2016-08-31 08:15:11"
